# Apply data-entry updates to the "relatorio_prontobaby_COMPLETO" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 - PACOTE PRÉ-OPERATÓRIO PEDIÁTRICO OTORRINO
$ws.Range("J2").Value = 6

# Row 3 - PACOTE PRÉ-OPERATÓRIO PEDIÁTRICO CIRURGIA GERAL
$ws.Range("H3").Value = 1
$ws.Range("K3").Value = 4
$ws.Range("M3").Value = 5

# Row 5 - ADENOIDECTOMIA PEDIÁTRICO
$ws.Range("J5").Value = 1

# Row 6 - fix procedure name spacing (remove space before dash)
$ws.Range("A6").Value = "AMIGDALECTOMIA- PEDIATRICO"

# Row 7 - AMIGDALECTOMIA COM ADENOIDECTOMIA - PEDIATRICO
$ws.Range("J7").Value = 3

# Row 11 - HERNIOPLASTIA UMBILICAL - PEDIATRICO
$ws.Range("K11").Value = 1

# Row 16 - POSTECTOMIA - PEDIATRICO
$ws.Range("K16").Value = 1
$ws.Range("M16").Value = 5

# Row 17 - TOTAL
$ws.Range("H17").Value = 1
$ws.Range("J17").Value = 10
$ws.Range("K17").Value = 6
$ws.Range("M17").Value = 10
